$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 14-27 (changed cells only) ---
# row 14
$ws.Range("D14").Value = 44435
$ws.Range("M14").Value = 115
$ws.Range("N14").Value = 18000
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 18000
$ws.Range("S14").Value = 1800
# row 15
$ws.Range("D15").Value = 44435
$ws.Range("M15").Value = 60
$ws.Range("N15").Value = 16000
$ws.Range("O15").Value = 16000
$ws.Range("P15").Value = 16000
$ws.Range("S15").Value = 1600
# row 16
$ws.Range("D16").Value = 44431
$ws.Range("M16").Value = 65
$ws.Range("N16").Value = 18000
$ws.Range("O16").Value = 18000
$ws.Range("P16").Value = 18000
$ws.Range("Q16").Value = "`$/bandeja 10 kilos"
$ws.Range("S16").Value = 1800
$ws.Range("T16").Value = 10
# row 17
$ws.Range("D17").Value = 44431
$ws.Range("L17").Value = "Segunda"
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 16000
$ws.Range("S17").Value = 1600
# row 18
$ws.Range("D18").Value = 44391
$ws.Range("L18").Value = "Primera"
$ws.Range("N18").Value = 17000
$ws.Range("O18").Value = 17000
$ws.Range("P18").Value = 17000
$ws.Range("S18").Value = 1700
# row 19
$ws.Range("D19").Value = 44391
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 45
$ws.Range("N19").Value = 15000
$ws.Range("O19").Value = 15000
$ws.Range("P19").Value = 15000
$ws.Range("S19").Value = 1500
# row 20
$ws.Range("D20").Value = 44315
$ws.Range("L20").Value = "Primera"
$ws.Range("N20").Value = 24000
$ws.Range("O20").Value = 24000
$ws.Range("P20").Value = 24000
$ws.Range("Q20").Value = "`$/caja 15 kilos granel"
$ws.Range("S20").Value = 1600
$ws.Range("T20").Value = 15
# row 21
$ws.Range("D21").Value = 44389
$ws.Range("M21").Value = 60
$ws.Range("N21").Value = 17000
$ws.Range("O21").Value = 17000
$ws.Range("P21").Value = 17000
$ws.Range("S21").Value = 1700
# row 22
$ws.Range("D22").Value = 44389
# row 23
$ws.Range("D23").Value = 44417
$ws.Range("M23").Value = 56
$ws.Range("N23").Value = 16000
$ws.Range("O23").Value = 16000
$ws.Range("P23").Value = 16000
$ws.Range("Q23").Value = "`$/bandeja 10 kilos"
$ws.Range("T23").Value = 10
# row 24
$ws.Range("D24").Value = 44417
$ws.Range("L24").Value = "Segunda"
$ws.Range("M24").Value = 60
$ws.Range("N24").Value = 14000
$ws.Range("O24").Value = 14000
$ws.Range("P24").Value = 14000
$ws.Range("S24").Value = 1400
# row 25
$ws.Range("D25").Value = 44420
$ws.Range("M25").Value = 54
$ws.Range("N25").Value = 18000
$ws.Range("O25").Value = 18000
$ws.Range("P25").Value = 18000
$ws.Range("S25").Value = 1800
# row 26
$ws.Range("D26").Value = 44420
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 50
$ws.Range("N26").Value = 15000
$ws.Range("O26").Value = 15000
$ws.Range("P26").Value = 15000
$ws.Range("S26").Value = 1500
# row 27
$ws.Range("D27").Value = 44433
$ws.Range("L27").Value = "Primera"
$ws.Range("N27").Value = 18000
$ws.Range("O27").Value = 18000
$ws.Range("P27").Value = 18000
$ws.Range("S27").Value = 1800

# --- Append new rows 28-32 ---
# row 28
$ws.Range("A28").Value = 3
$ws.Range("B28").Value = "Femacal de La Calera"
$ws.Range("C28").Value = "Coquimbo"
$ws.Range("D28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D28").Value = 44319
$ws.Range("E28").Value = 5
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100108
$ws.Range("H28").Value = "Tropicales y subtropicales"
$ws.Range("I28").Value = 100108004
$ws.Range("J28").Value = "Papaya"
$ws.Range("K28").Value = "Cultivar IV Región"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 60
$ws.Range("N28").Value = 24000
$ws.Range("O28").Value = 24000
$ws.Range("P28").Value = 24000
$ws.Range("Q28").Value = "`$/caja 15 kilos granel"
$ws.Range("R28").Value = "Provincia del Elquí"
$ws.Range("S28").Value = 1600
$ws.Range("T28").Value = 15
# row 29
$ws.Range("A29").Value = 3
$ws.Range("B29").Value = "Femacal de La Calera"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D29").Value = 44370
$ws.Range("E29").Value = 5
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100108
$ws.Range("H29").Value = "Tropicales y subtropicales"
$ws.Range("I29").Value = 100108004
$ws.Range("J29").Value = "Papaya"
$ws.Range("K29").Value = "Cultivar IV Región"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 50
$ws.Range("N29").Value = 17000
$ws.Range("O29").Value = 17000
$ws.Range("P29").Value = 17000
$ws.Range("Q29").Value = "`$/bandeja 10 kilos"
$ws.Range("R29").Value = "Provincia del Elquí"
$ws.Range("S29").Value = 1700
$ws.Range("T29").Value = 10
# row 30
$ws.Range("A30").Value = 3
$ws.Range("B30").Value = "Femacal de La Calera"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D30").Value = 44382
$ws.Range("E30").Value = 5
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100108
$ws.Range("H30").Value = "Tropicales y subtropicales"
$ws.Range("I30").Value = 100108004
$ws.Range("J30").Value = "Papaya"
$ws.Range("K30").Value = "Cultivar IV Región"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 58
$ws.Range("N30").Value = 17000
$ws.Range("O30").Value = 17000
$ws.Range("P30").Value = 17000
$ws.Range("Q30").Value = "`$/bandeja 10 kilos"
$ws.Range("R30").Value = "Provincia del Elquí"
$ws.Range("S30").Value = 1700
$ws.Range("T30").Value = 10
# row 31
$ws.Range("A31").Value = 3
$ws.Range("B31").Value = "Femacal de La Calera"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D31").Value = 44398
$ws.Range("E31").Value = 5
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100108
$ws.Range("H31").Value = "Tropicales y subtropicales"
$ws.Range("I31").Value = 100108004
$ws.Range("J31").Value = "Papaya"
$ws.Range("K31").Value = "Cultivar IV Región"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 60
$ws.Range("N31").Value = 17000
$ws.Range("O31").Value = 17000
$ws.Range("P31").Value = 17000
$ws.Range("Q31").Value = "`$/bandeja 10 kilos"
$ws.Range("R31").Value = "Provincia del Elquí"
$ws.Range("S31").Value = 1700
$ws.Range("T31").Value = 10
# row 32
$ws.Range("A32").Value = 3
$ws.Range("B32").Value = "Femacal de La Calera"
$ws.Range("C32").Value = "Coquimbo"
$ws.Range("D32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D32").Value = 44398
$ws.Range("E32").Value = 5
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100108
$ws.Range("H32").Value = "Tropicales y subtropicales"
$ws.Range("I32").Value = 100108004
$ws.Range("J32").Value = "Papaya"
$ws.Range("K32").Value = "Cultivar IV Región"
$ws.Range("L32").Value = "Segunda"
$ws.Range("M32").Value = 50
$ws.Range("N32").Value = 15000
$ws.Range("O32").Value = 15000
$ws.Range("P32").Value = 15000
$ws.Range("Q32").Value = "`$/bandeja 10 kilos"
$ws.Range("R32").Value = "Provincia del Elquí"
$ws.Range("S32").Value = 1500
$ws.Range("T32").Value = 10

# dimension will auto-update to A1:T32 on save
